$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.636.61"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "1.963.64"
$ws.Range("E3").Value = "  +1.76%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("E6").Value = "  +1.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.20%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +5.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0794"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.95%  "

$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("E12").Value = "  +7.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.841"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.04%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.99%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.252.76"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.35%  "

$ws.Range("D17").Value = "1.965.11"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("D18").Value = "36.567.19"
$ws.Range("E18").Value = "  +0.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.34%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.36%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  +7.50%  "

$ws.Range("E25").Value = "  +4.04%  "

$ws.Range("E26").Value = "  +11.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("E30").Value = "  +11.79%  "

$ws.Range("E31").Value = "  +2.97%  "

$ws.Range("E32").Value = "  +6.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0619"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +20.68%  "

$ws.Range("E36").Value = "  +7.63%  "

$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.98%  "

$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("E42").Value = "  +4.38%  "

$ws.Range("E43").Value = "  +2.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("D45").Value = "1.367.41"
$ws.Range("E45").Value = "  +3.68%  "

$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.44%  "

$ws.Range("E51").Value = "  +6.35%  "
